# Updates the cryptos list worksheet with the latest scraped price/volume
# figures (and a few re-sorted rows) per the GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay plain text (coinranking prices use
# dotted thousands-separators, e.g. '41.215.80', and some look like plain
# numbers, e.g. 0.616 -- force text format so Excel does not reinterpret
# them as numeric values, then drop the format override again so the
# cell keeps its original (default) style.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
$ws.Range('D2').Value = '41.215.80'
$ws.Range('E2').Value = '  -1.62%  '
# Row 3
$ws.Range('D3').Value = '2.182.21'
$ws.Range('E3').Value = '  -2.20%  '
# Row 4
$ws.Range('E4').Value = '  -0.10%  '
# Row 5
Set-TextValue $ws.Range('D5') '249.70'
$ws.Range('E5').Value = '  +0.18%  '
# Row 6
Set-TextValue $ws.Range('D6') '0.616'
$ws.Range('E6').Value = '  -2.34%  '
# Row 7
Set-TextValue $ws.Range('D7') '66.45'
$ws.Range('E7').Value = '  -7.39%  '
# Row 8
$ws.Range('E8').Value = '  -0.04%  '
# Row 9
Set-TextValue $ws.Range('D9') '0.575'
$ws.Range('E9').Value = '  -3.83%  '
# Row 10
Set-TextValue $ws.Range('D10') '58.77'
$ws.Range('E10').Value = '  +1.23%  '
# Row 11
Set-TextValue $ws.Range('D11') '36.25'
$ws.Range('E11').Value = '  -10.73%  '
# Row 12
Set-TextValue $ws.Range('D12') '0.0931'
$ws.Range('E12').Value = '  -4.23%  '
# Row 13
Set-TextValue $ws.Range('D13') '0.103'
$ws.Range('E13').Value = '  -2.62%  '
# Row 14
Set-TextValue $ws.Range('D14') '6.89'
$ws.Range('E14').Value = '  -3.69%  '
# Row 15
$ws.Range('D15').Value = '2.510.46'
$ws.Range('E15').Value = '  -2.23%  '
# Row 16
Set-TextValue $ws.Range('D16') '14.37'
$ws.Range('E16').Value = '  -4.08%  '
# Row 17
$ws.Range('E17').Value = '  -1.20%  '
# Row 18
$ws.Range('D18').Value = '2.171.23'
$ws.Range('E18').Value = '  -2.43%  '
# Row 19
$ws.Range('D19').Value = '41.166.56'
$ws.Range('E19').Value = '  -1.89%  '
# Row 20
$ws.Range('E20').Value = '  -2.53%  '
# Row 21
Set-TextValue $ws.Range('D21') '71.70'
$ws.Range('E21').Value = '  -1.95%  '
# Row 22
$ws.Range('E22').Value = '  -2.54%  '
# Row 23
Set-TextValue $ws.Range('D23') '230.47'
$ws.Range('E23').Value = '  -2.26%  '
# Row 24
Set-TextValue $ws.Range('D24') '2.03'
$ws.Range('E24').Value = '  -4.34%  '
# Row 25
Set-TextValue $ws.Range('D25') '3.81'
$ws.Range('E25').Value = '  -5.48%  '
# Row 26
$ws.Range('E26').Value = '  +0.06%  '
# Row 27
Set-TextValue $ws.Range('D27') '11.37'
$ws.Range('E27').Value = '  +6.22%  '
# Row 28
$ws.Range('E28').Value = '  -4.86%  '
# Row 29
$ws.Range('E29').Value = '  -3.91%  '
# Row 30
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range('D30') '2.12'
$ws.Range('E30').Value = '  -3.28%  '
# Row 31
$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range('D31') '168.23'
$ws.Range('E31').Value = '  -1.99%  '
# Row 32
Set-TextValue $ws.Range('D32') '20.25'
$ws.Range('E32').Value = '  -2.40%  '
# Row 33
$ws.Range('E33').Value = '  -1.54%  '
# Row 34
Set-TextValue $ws.Range('D34') '5.77'
$ws.Range('E34').Value = '  +4.08%  '
# Row 35
Set-TextValue $ws.Range('D35') '0.0743'
$ws.Range('E35').Value = '  +1.12%  '
# Row 36
$ws.Range('E36').Value = '  -3.03%  '
# Row 37
$ws.Range('E37').Value = '  -4.33%  '
# Row 38
Set-TextValue $ws.Range('D38') '3.97'
$ws.Range('E38').Value = '  -0.96%  '
# Row 39
Set-TextValue $ws.Range('D39') '24.60'
$ws.Range('E39').Value = '  -7.26%  '
# Row 40
$ws.Range('E40').Value = '  +3.69%  '
# Row 41
$ws.Range('E41').Value = '  -3.41%  '
# Row 42
Set-TextValue $ws.Range('D42') '5.31'
$ws.Range('E42').Value = '  +7.15%  '
# Row 43
$ws.Range('E43').Value = '  -8.48%  '
# Row 44
$ws.Range('B44').Value = 'MultiversX'
$ws.Range('C44').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-TextValue $ws.Range('D44') '61.30'
$ws.Range('E44').Value = '  -7.87%  '
# Row 45
$ws.Range('B45').Value = 'Celestia'
$ws.Range('C45').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue $ws.Range('D45') '11.43'
$ws.Range('E45').Value = '  -5.49%  '
# Row 46
$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range('D46') '0.190'
$ws.Range('E46').Value = '  -7.52%  '
# Row 47
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D47') '8.51'
$ws.Range('E47').Value = '  -3.18%  '
# Row 48
Set-TextValue $ws.Range('D48') '0.0999'
$ws.Range('E48').Value = '  -1.57%  '
# Row 49
$ws.Range('E49').Value = '  -0.08%  '
# Row 50
$ws.Range('E50').Value = '  -1.79%  '
# Row 51
$ws.Range('E51').Value = '  -3.85%  '
